$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.104.90'
$ws.Range("E2").Value = '  -3.29%  '
$ws.Range("D3").Value = '2.367.74'
$ws.Range("E3").Value = '  -3.52%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '501.05'
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.68'
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -2.46%  '
$ws.Range("D9").Value = '2.371.63'
$ws.Range("E9").Value = '  -3.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0980'
$ws.Range("E10").Value = '  +0.37%  '
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.78'
$ws.Range("E12").Value = '  +3.69%  '
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").Value = '2.787.06'
$ws.Range("E14").Value = '  -3.39%  '
$ws.Range("D15").Value = '56.008.04'
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.44'
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = '2.397.67'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.00'
$ws.Range("E19").Value = '  -3.00%  '
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '306.57'
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.29'
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.89'
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.367'
$ws.Range("E26").Value = '  -3.70%  '
$ws.Range("E27").Value = '  -5.83%  '
$ws.Range("E28").Value = '  -4.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.72'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("E31").Value = '  -3.50%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.997'
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("E34").Value = '  -7.00%  '
$ws.Range("E35").Value = '  -5.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.60'
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("E37").Value = '  -5.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.74'
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.788'
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("E41").Value = '  -5.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '129.86'
$ws.Range("E42").Value = '  -4.92%  '
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("E44").Value = '  -6.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.564'
$ws.Range("E45").Value = '  -1.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0902'
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '239.85'
$ws.Range("E47").Value = '  -6.61%  '
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("E49").Value = '  -3.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.05'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.949'
$ws.Range("E51").Value = '  -0.72%  '

$resetCells = @("D5","D6","D10","D12","D16","D19","D21","D22","D24","D25","D26","D29","D33","D36","D38","D40","D42","D45","D46","D47","D50","D51")
foreach ($addr in $resetCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
